$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert the two new "Thẩm quyền" columns into the table (cols K and, after
#    the first insert shifts things over, the slot that will become N), then
#    delete the now-orphaned old "Kết quả" column that gets pushed out to the
#    far right.
# ---------------------------------------------------------------------------
$ws.Range("K1").EntireColumn.Insert()
$ws.Range("N1").EntireColumn.Insert()
$ws.Range("P1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2) Move the footer note from I21 down to J22 (row 21 -> 22, col I -> J) -
#    this does NOT shift any other rows (row 26 stays put).
# ---------------------------------------------------------------------------
$noteValue = $ws.Range("I21").Value()
$ws.Range("D5").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("J22").Value = $noteValue
$ws.Range("I21").Clear()

# ---------------------------------------------------------------------------
# 3) Header row 16: fill in the two new column headers.
# ---------------------------------------------------------------------------
$ws.Range("K16").Value = "Thẩm quyền GQKN lần I"
$ws.Range("N16").Value = "Thẩm quyền GQKN lần II"

# Row 17 is a plain 1..15 index row - fill the two new cells and extend.
$ws.Range("K17").Value = 11
$ws.Range("N17").Value = 14
$ws.Range("O17").Value = 15

# ---------------------------------------------------------------------------
# 4) Row 18 data: K18/N18 should be blank with the plain bordered style (5),
#    not inherit the date-format style that Insert copied from J/M.
# ---------------------------------------------------------------------------
$ws.Range("B18").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("N18").PasteSpecial(-4122)
$ws.Range("K18").ClearContents()
$ws.Range("N18").ClearContents()
